$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two remaining "NUMERO DE RADICACION" records with corrected values
$ws.Range("B2").Value2 = "11001310300520210054100"
$ws.Range("B3").Value2 = "11001310300420170072500"

# The old run had produced extra (erroneous) records in B4:B7 -- remove them
# entirely so the sheet only keeps the header + the two corrected rows.
$ws.Range("B4:B7").Clear()

# Reflect the new selection the author ended up with after trimming the table
$ws.Range("B4").Select() | Out-Null
